# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> (old, new) value map is identical on both sheets, but the row
# numbers differ because 全部类型 interleaves rows from multiple sheets.
$exhibitRows = @{
    2  = 270
    14 = 2585
    16 = 1406
    17 = 5054
    21 = 1928
    23 = 3357
    28 = 140
    29 = 6
    30 = 329
    32 = 2152
    33 = 4
    36 = 793
}

$allTypeRows = @{
    7  = 270
    20 = 2585
    21 = 1406
    26 = 5054
    30 = 1928
    32 = 3357
    40 = 140
    41 = 6
    42 = 329
    44 = 2152
    45 = 4
    48 = 793
}

foreach ($row in $exhibitRows.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitRows[$row]
}

foreach ($row in $allTypeRows.Keys) {
    $wsAll.Range("F$row").Value = $allTypeRows[$row]
}
